$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update the Maximo (C2) metric ---
$resumen = $wb.Worksheets.Item("Resumen")
$resumen.Range("C2").Value = 548.0019785922802

# --- Sheet "Solucion": randomized Pedido order + recomputed Salida (Zona-route) assignment ---
$solucion = $wb.Worksheets.Item("Solucion")
$data = @(
  @("Pedido_4","S001"),
  @("Pedido_27","S021"),
  @("Pedido_45","S041"),
  @("Pedido_48","S051"),
  @("Pedido_54","S031"),
  @("Pedido_8","S011"),
  @("Pedido_10","S042"),
  @("Pedido_37","S002"),
  @("Pedido_56","S022"),
  @("Pedido_6","S052"),
  @("Pedido_22","S032"),
  @("Pedido_20","S012"),
  @("Pedido_35","S043"),
  @("Pedido_16","S023"),
  @("Pedido_30","S003"),
  @("Pedido_49","S053"),
  @("Pedido_38","S033"),
  @("Pedido_52","S013"),
  @("Pedido_26","S044"),
  @("Pedido_39","S024"),
  @("Pedido_9","S004"),
  @("Pedido_58","S034"),
  @("Pedido_11","S054"),
  @("Pedido_2","S014"),
  @("Pedido_25","S025"),
  @("Pedido_5","S045"),
  @("Pedido_28","S005"),
  @("Pedido_3","S035"),
  @("Pedido_41","S055"),
  @("Pedido_19","S015"),
  @("Pedido_60","S026"),
  @("Pedido_29","S046"),
  @("Pedido_23","S006"),
  @("Pedido_57","S036"),
  @("Pedido_7","S027"),
  @("Pedido_18","S056"),
  @("Pedido_32","S016"),
  @("Pedido_13","S037"),
  @("Pedido_36","S047"),
  @("Pedido_42","S007"),
  @("Pedido_15","S057"),
  @("Pedido_59","S028"),
  @("Pedido_17","S017"),
  @("Pedido_44","S038"),
  @("Pedido_14","S048"),
  @("Pedido_55","S008"),
  @("Pedido_21","S029"),
  @("Pedido_43","S058"),
  @("Pedido_51","S018"),
  @("Pedido_24","S049"),
  @("Pedido_1","S009"),
  @("Pedido_33","S039"),
  @("Pedido_47","S059"),
  @("Pedido_40","S019"),
  @("Pedido_50","S030"),
  @("Pedido_34","S050"),
  @("Pedido_31","S010"),
  @("Pedido_46","S040"),
  @("Pedido_53","S060"),
  @("Pedido_12","S020")
)
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $solucion.Cells.Item($row, 1).Value = $data[$i][0]
  $solucion.Cells.Item($row, 2).Value = $data[$i][1]
}

# --- Sheet "Metricas": update per-zone time metrics ---
$metricas = $wb.Worksheets.Item("Metricas")
$metricas.Range("B2").Value = 548.0019785922802
$metricas.Range("B3").Value = 529.9498540382743
$metricas.Range("B4").Value = 544.8734944318305
